$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts existing E/F -> F/G), matching the
# "location" column being duplicated with an alternate (en-dash) value.
$ws.Columns("E:E").Insert()

# Header for the new column (same text as D1 "location")
$ws.Range("E1").Value = "location"

# New location string - note this uses an EN DASH (U+2013), not a hyphen,
# between "qml" and "21".
$ws.Range("E2").Value = "../../QML/OggettiEditDash/AddNewForm.qml " + [char]0x2013 + " 21"

# Match the formatting of the neighboring "location" cell (D2) for the new
# E2 cell, which picks up its own distinct cell style.
$ws.Range("E2").Style = $ws.Range("D2").Style

# Update the active selection to match the edited workbook's cursor position.
$null = $ws.Range("E8").Select()
